$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowPair($ws, $r1, $r2) {
    foreach ($col in @("B","D","E","F","G")) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

$pairs = @(
    @(192,193),
    @(227,228),
    @(322,323),
    @(364,365),
    @(366,367),
    @(372,373),
    @(375,376),
    @(380,381),
    @(442,443),
    @(463,464),
    @(572,573)
)

foreach ($p in $pairs) {
    Swap-RowPair $ws $p[0] $p[1]
}
